$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. The opening paragraph of the use-case diagram section had its text
#    split across two runs ("...principales fo" / "nctionnalités...")
#    with a _GoBack bookmark sandwiched in between. Re-typing the full
#    sentence over itself merges it back into a single run and drops
#    that stray bookmark.
# ---------------------------------------------------------------------
$sentence = "Ce diagramme montre les principales fonctionnalités du système et les acteurs qui interagissent avec lui."
$d.Content.Find.Execute($sentence, $true, $false, $false, $false, $false, $true, 1, $false, $sentence, 2) | Out-Null

# ---------------------------------------------------------------------
# 2. The use-case list had a "Choix des sièges" bullet that needs to go
#    away entirely; the rest of the list ("Saisie des informations
#    personnelles", "Paiement", "Confirmation de la réservation") shifts
#    up to fill the gap.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Choix*") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 3. Re-insert a _GoBack bookmark right after the word "paiement" in the
#    "Passerelle de paiement" bullet (component diagram section). A
#    directly-collapsed range at that exact paragraph-ending offset is
#    unreliable, so nudge past the end with a throw-away character,
#    plant the bookmark just before it, then remove the character again.
# ---------------------------------------------------------------------
$hit = $d.Content
$ok = $hit.Find.Execute("Passerelle de paiement", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($ok) {
    $tail = $hit.End
    $anchor = $d.Range($tail, $tail)
    $anchor.InsertAfter("#")
    $bm = $d.Range($tail, $tail)
    $d.Bookmarks.Add("_GoBack", $bm)
    $d.Range($tail, $tail + 1).Delete()
}
